$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 14 data: Audit_time / Date / Time when the audit was recorded
$ws.Cells.Item(14, 2).Value = 12
$ws.Cells.Item(14, 3).Value = "Audit_time"
$ws.Cells.Item(14, 4).Value = "Date"
$ws.Cells.Item(14, 5).Value = "Time when the audit was recorded"

# Match bold style used by other rows in column B (B3:B13)
$ws.Cells.Item(14, 2).Font.Bold = $true

# Update the active selection to E14 as in the diff
$ws.Range("E14").Select()
